$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 35: new transaction - 3 botellones on 2023-09-19, amount -159
$ws.Cells.Item(35, 2).Value = 45188          # B35 fecha
$ws.Cells.Item(35, 3).Value = "3 botellones" # C35 descripcion
$ws.Cells.Item(35, 4).Value = -159           # D35 importe

# Row 36: new transaction - 3 botellones on 2023-09-22, amount -159
$ws.Cells.Item(36, 2).Value = 45191          # B36 fecha
$ws.Cells.Item(36, 3).Value = "3 botellones" # C36 descripcion
$ws.Cells.Item(36, 4).Value = -159           # D36 importe

# Update the active selection to match the author's final cursor position
$ws.Range("D37").Select()
